$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1792.1765
$ws.Range("I2").Value = 1154
$ws.Range("J2").Value = 2703.8572
$ws.Range("K2").Value = 1154
$ws.Range("L2").Value = 2703.8572
$ws.Range("M2").Value = -1041
$ws.Range("N2").Value = -2929.8572
$ws.Range("H18").Value = 936
$ws.Range("I18").Value = 676.25
$ws.Range("K18").Value = 676.25
$ws.Range("M18").Value = -392.25
$ws.Range("H19").Value = 1244.5
$ws.Range("I19").Value = 1405.5
$ws.Range("J19").Value = 1180.1
$ws.Range("K19").Value = 1405.5
$ws.Range("L19").Value = 1180.1
$ws.Range("M19").Value = -1230.5
$ws.Range("N19").Value = -1530.1
$ws.Range("H88").Value = 1214.381
$ws.Range("J88").Value = 1833.2
$ws.Range("L88").Value = 1833.2
$ws.Range("N88").Value = -2645.2
$ws.Range("H91").Value = 1214.381
$ws.Range("J91").Value = 1833.2
$ws.Range("L91").Value = 1833.2
$ws.Range("N91").Value = -4641.2
$ws.Range("H103").Value = 1766.3334
$ws.Range("I103").Value = 3499
$ws.Range("J103").Value = 900
$ws.Range("K103").Value = 10497
$ws.Range("L103").Value = 2700
$ws.Range("M103").Value = -9911
$ws.Range("N103").Value = -3872
$ws.Range("H132").Value = 2852.5
$ws.Range("I132").Value = 2324.348
$ws.Range("K132").Value = 6973.044
$ws.Range("M132").Value = -4443.044
$ws.Range("H133").Value = 124995.5
$ws.Range("J133").Value = 124995.5
$ws.Range("L133").Value = 124995.5
$ws.Range("N133").Value = -135115.5
$ws.Range("H136").Value = 125779.5
$ws.Range("J136").Value = 125779.5
$ws.Range("L136").Value = 125779.5
$ws.Range("N136").Value = -135979.5
$ws.Range("H141").Value = 6382.294
$ws.Range("I141").Value = 4049.9
$ws.Range("J141").Value = 9714.286
$ws.Range("K141").Value = 12149.7
$ws.Range("L141").Value = 29142.858
$ws.Range("M141").Value = -6969.700000000001
$ws.Range("N141").Value = -39502.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 267.57144
$ws.Range("J5").Value = 136.33333
$ws.Range("L5").Value = 136.33333
$ws.Range("N5").Value = -360.33333
$ws.Range("H45").Value = 2676.3076
$ws.Range("I45").Value = 1499.5
$ws.Range("J45").Value = 4559.2
$ws.Range("K45").Value = 1499.5
$ws.Range("L45").Value = 4559.2
$ws.Range("M45").Value = -1122.5
$ws.Range("N45").Value = -5313.2
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 60000
$ws.Range("M65").Value = -56880
$ws.Range("H74").Value = 3335.0527
$ws.Range("I74").Value = 2642.5557
$ws.Range("K74").Value = 2642.5557
$ws.Range("M74").Value = -1768.5557
$ws.Range("H77").Value = 3335.0527
$ws.Range("I77").Value = 2642.5557
$ws.Range("K77").Value = 13212.7785
$ws.Range("M77").Value = -8844.7785
$ws.Range("H122").Value = 3316.55
$ws.Range("I122").Value = 2435.1667
$ws.Range("J122").Value = 11249
$ws.Range("K122").Value = 7305.500100000001
$ws.Range("L122").Value = 33747
$ws.Range("M122").Value = -4855.500100000001
$ws.Range("N122").Value = -38647
$ws.Range("H132").Value = 2785.8235
$ws.Range("I132").Value = 2324.7878
$ws.Range("K132").Value = 6974.3634
$ws.Range("M132").Value = -4444.3634

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 267.57144
$ws.Range("J4").Value = 136.33333
$ws.Range("L4").Value = 136.33333
$ws.Range("N4").Value = -366.33333
$ws.Range("H22").Value = 1066.1
$ws.Range("I22").Value = 926.8333
$ws.Range("K22").Value = 926.8333
$ws.Range("M22").Value = -753.8333
$ws.Range("H45").Value = 44000
$ws.Range("J45").Value = 44000
$ws.Range("L45").Value = 44000
$ws.Range("N45").Value = -45616
$ws.Range("H59").Value = 94779.89999999999
$ws.Range("J59").Value = 96422.11
$ws.Range("L59").Value = 96422.11
$ws.Range("N59").Value = -98116.11
$ws.Range("H138").Value = 110000
$ws.Range("I138").Value = 60000
$ws.Range("J138").Value = 135000
$ws.Range("K138").Value = 60000
$ws.Range("L138").Value = 135000
$ws.Range("M138").Value = -54860
$ws.Range("N138").Value = -145280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6882.44
$ws.Range("I31").Value = 2433.7097
$ws.Range("J31").Value = 14140.895
$ws.Range("K31").Value = 2433.7097
$ws.Range("L31").Value = 14140.895
$ws.Range("M31").Value = -2138.7097
$ws.Range("N31").Value = -14730.895
$ws.Range("H34").Value = 6882.44
$ws.Range("I34").Value = 2433.7097
$ws.Range("J34").Value = 14140.895
$ws.Range("K34").Value = 2433.7097
$ws.Range("L34").Value = 14140.895
$ws.Range("M34").Value = -2231.7097
$ws.Range("N34").Value = -14544.895
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H135").Value = 93889
$ws.Range("J135").Value = 93889
$ws.Range("L135").Value = 93889
$ws.Range("N135").Value = -104029
$ws.Range("H140").Value = 89665
$ws.Range("J140").Value = 89665
$ws.Range("L140").Value = 89665
$ws.Range("N140").Value = -100025

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 820.61536
$ws.Range("J86").Value = 891.1429000000001
$ws.Range("L86").Value = 2673.4287
$ws.Range("N86").Value = -5045.4287
$ws.Range("H89").Value = 820.61536
$ws.Range("J89").Value = 891.1429000000001
$ws.Range("L89").Value = 8020.2861
$ws.Range("N89").Value = -19876.2861
$ws.Range("H131").Value = 4642.375
$ws.Range("J131").Value = 4786
$ws.Range("L131").Value = 14358
$ws.Range("N131").Value = -24438
$ws.Range("H132").Value = 3608.2666
$ws.Range("I132").Value = 3608.2666
$ws.Range("K132").Value = 32474.3994
$ws.Range("M132").Value = -29944.3994
$ws.Range("H141").Value = 6160.294
$ws.Range("I141").Value = 4963.3335
$ws.Range("K141").Value = 14890.0005
$ws.Range("M141").Value = -9710.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 577
$ws.Range("I97").Value = 430.36365
$ws.Range("J97").Value = 980.25
$ws.Range("K97").Value = 430.36365
$ws.Range("L97").Value = 980.25
$ws.Range("M97").Value = 65.63634999999999
$ws.Range("N97").Value = -1972.25
$ws.Range("H114").Value = 250000
$ws.Range("J114").Value = 250000
$ws.Range("L114").Value = 250000
$ws.Range("N114").Value = -258678
$ws.Range("H118").Value = 20577.75
$ws.Range("J118").Value = 20577.75
$ws.Range("L118").Value = 20577.75
$ws.Range("N118").Value = -23891.75
$ws.Range("H122").Value = 12278.8
$ws.Range("J122").Value = 14079.6
$ws.Range("L122").Value = 42238.8
$ws.Range("N122").Value = -47138.8
$ws.Range("H132").Value = 3674.9744
$ws.Range("I132").Value = 2942.743
$ws.Range("J132").Value = 10082
$ws.Range("K132").Value = 8828.228999999999
$ws.Range("L132").Value = 30246
$ws.Range("M132").Value = -6298.228999999999
$ws.Range("N132").Value = -35306
$ws.Range("H135").Value = 120000
$ws.Range("J135").Value = 120000
$ws.Range("L135").Value = 120000
$ws.Range("N135").Value = -130140
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1059.7391
$ws.Range("I55").Value = 191.90909
$ws.Range("K55").Value = 191.90909
$ws.Range("M55").Value = -18.90908999999999
$ws.Range("H122").Value = 5476.5713
$ws.Range("I122").Value = 4318.1177
$ws.Range("J122").Value = 10400
$ws.Range("K122").Value = 12954.3531
$ws.Range("L122").Value = 31200
$ws.Range("M122").Value = -10504.3531
$ws.Range("N122").Value = -36100
$ws.Range("H132").Value = 8674.15
$ws.Range("I132").Value = 9483.286
$ws.Range("K132").Value = 28449.858
$ws.Range("M132").Value = -25919.858
$ws.Range("H136").Value = 8343.189
$ws.Range("I136").Value = 4403.909
$ws.Range("K136").Value = 13211.727
$ws.Range("M136").Value = -10661.727
$ws.Range("H140").Value = 89999.5
$ws.Range("J140").Value = 89999.5
$ws.Range("L140").Value = 89999.5
$ws.Range("N140").Value = -100359.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11609.889
$ws.Range("J62").Value = 14697.8
$ws.Range("L62").Value = 14697.8
$ws.Range("N62").Value = -15945.8
$ws.Range("H65").Value = 11609.889
$ws.Range("J65").Value = 14697.8
$ws.Range("L65").Value = 73489
$ws.Range("N65").Value = -79729
$ws.Range("H95").Value = 153943
$ws.Range("J95").Value = 153943
$ws.Range("L95").Value = 153943
$ws.Range("N95").Value = -159435
$ws.Range("H126").Value = 3797.7896
$ws.Range("I126").Value = 3797.7896
$ws.Range("K126").Value = 11393.3688
$ws.Range("M126").Value = -8923.3688
$ws.Range("H132").Value = 2859.853
$ws.Range("I132").Value = 1628.862
$ws.Range("K132").Value = 4886.586
$ws.Range("M132").Value = -2356.586
$ws.Range("H136").Value = 3416.647
$ws.Range("I136").Value = 3055.1785
$ws.Range("J136").Value = 5103.5
$ws.Range("K136").Value = 9165.5355
$ws.Range("L136").Value = 15310.5
$ws.Range("M136").Value = -6615.5355
$ws.Range("N136").Value = -20410.5

Write-Output "Applied 243 cell updates across 8 sheets"
